# Actualizacion automatica 2025-09-01 08:30:07
#
# This script rolls the monthly sales tracker forward by one month:
#   - "VENTAS POR GRUPO": the per-product-group figures that had been
#     accumulated for the closing month are cleared back to 0 (and the
#     "x de 55" coverage counters on row 57 reset to "0 de 55"), since
#     that month's column in "VENTA MENSUAL" is about to roll off.
#   - "VENTA MENSUAL": the monthly columns C:F (previously
#     mayo/junio/julio/agosto) shift one month to the left
#     (junio/julio/agosto/septiembre), carrying each client's totals
#     forward a column and opening a fresh (all-zero) column for the
#     newest month.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet1 "VENTAS POR GRUPO": zero out sales figures for specific clients/products ---
$ws1.Range("L5").Value = 0
$ws1.Range("D6").Value = 0
$ws1.Range("K6").Value = 0
$ws1.Range("L6").Value = 0
$ws1.Range("M6").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("O14").Value = 0
$ws1.Range("C22").Value = 0
$ws1.Range("D22").Value = 0
$ws1.Range("M22").Value = 0
$ws1.Range("C25").Value = 0
$ws1.Range("D26").Value = 0
$ws1.Range("D28").Value = 0
$ws1.Range("K29").Value = 0
$ws1.Range("Q29").Value = 0
$ws1.Range("M30").Value = 0
$ws1.Range("I31").Value = 0
$ws1.Range("M31").Value = 0
$ws1.Range("C33").Value = 0
$ws1.Range("C38").Value = 0
$ws1.Range("D38").Value = 0
$ws1.Range("E38").Value = 0
$ws1.Range("L38").Value = 0
$ws1.Range("M38").Value = 0
$ws1.Range("D41").Value = 0
$ws1.Range("K41").Value = 0
$ws1.Range("M41").Value = 0
$ws1.Range("L43").Value = 0
$ws1.Range("D44").Value = 0
$ws1.Range("H44").Value = 0
$ws1.Range("L44").Value = 0
$ws1.Range("M44").Value = 0
$ws1.Range("C45").Value = 0
$ws1.Range("D45").Value = 0
$ws1.Range("H45").Value = 0
$ws1.Range("I45").Value = 0
$ws1.Range("L46").Value = 0
$ws1.Range("M46").Value = 0
$ws1.Range("H47").Value = 0
$ws1.Range("I47").Value = 0
$ws1.Range("M47").Value = 0
$ws1.Range("M49").Value = 0
$ws1.Range("R49").Value = 0
$ws1.Range("I54").Value = 0

# --- Sheet1 row 57: refresh "x de 55" coverage counters to "0 de 55" ---
$ws1.Range("C57").Value = "0 de 55"
$ws1.Range("D57").Value = "0 de 55"
$ws1.Range("E57").Value = "0 de 55"
$ws1.Range("H57").Value = "0 de 55"
$ws1.Range("I57").Value = "0 de 55"
$ws1.Range("K57").Value = "0 de 55"
$ws1.Range("L57").Value = "0 de 55"
$ws1.Range("M57").Value = "0 de 55"
$ws1.Range("O57").Value = "0 de 55"
$ws1.Range("Q57").Value = "0 de 55"
$ws1.Range("R57").Value = "0 de 55"

# --- Sheet2 "VENTA MENSUAL": roll the month headers forward ---
$ws2.Range("C1").Value = "junio"
$ws2.Range("D1").Value = "julio"
$ws2.Range("E1").Value = "agosto"
$ws2.Range("F1").Value = "septiembre"

# Column width tweaks that came with the new month column.
# (raw OOXML <col width> = ColumnWidth + 0.8333333333333334 under this engine)
$ws2.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 15.166666666666666

# Shift each client's monthly figures one column to the left:
# new C = old D, new D = old E, new E = old F, new F = 0 (brand-new month).
# Read with .Value2 (reliable read accessor here) before writing with .Value.
for ($r = 2; $r -le 57; $r++) {
    $oldD = $ws2.Cells.Item($r, 4).Value2
    $oldE = $ws2.Cells.Item($r, 5).Value2
    $oldF = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($r, 3).Value = $oldD
    $ws2.Cells.Item($r, 4).Value = $oldE
    $ws2.Cells.Item($r, 5).Value = $oldF
    $ws2.Cells.Item($r, 6).Value = 0
}
